# Auto-committed on 2022/01/20 週四
# Update the numbered-list punctuation (Chinese full-width period "." -> colon ":")
# in the "DBD" sheet's remark column (G), rows 9 and 10, and move the
# active selection/viewport the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")   # "DBD" is the sheet shown/selected in the source file

# Row 9, column G: "1.排除商品別 / 2.排除部門別 / 3.是否排除15日薪非業績人員 /
# 4.業績追回時通知員工代碼(email) / R.業績重算設定"
# -> change each list-marker's full stop to a colon.
$ws.Range("G9").Value = "1:排除商品別`n2:排除部門別`n3:是否排除15日薪非業績人員`n4:業績追回時通知員工代碼(email)`nR:業績重算設定"

# Row 10, column G: "條件記號1為1、2、3時 / 1.業績全部 / 2.換算業績、業務報酬 /
# 3.介紹獎金 / 4.加碼獎勵津貼 / 5.協辦獎金 / 其他為空白1位"
# -> change each list-marker's full stop to a colon.
$ws.Range("G10").Value = "條件記號1為1、2、3時`n1:業績全部`n2:換算業績、業務報酬`n3:介紹獎金`n4:加碼獎勵津貼`n5:協辦獎金`n其他為空白1位"

# Move the viewport / selection on the DBD sheet (topLeftCell A7 -> A10,
# selection H11 -> G14).
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G14").Select() | Out-Null
